$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Existing-cell tweaks
# ---------------------------------------------------------------------------

# M16: learning rate value change
$ws.Range("M16").Value = 0.0105

# Column L (12) width change (raw OOXML width 6 -> 9). ColumnWidth 8.15
# (Excel "characters" units) round-trips to a stored width of exactly 9.
$ws.Columns.Item(12).ColumnWidth = 8.15

# ---------------------------------------------------------------------------
# 2) New data rows (24-26, 30-32, 35-38, 42-44)
#
# The *order* in which brand-new text values are first written controls the
# order they are appended to the shared-string table, so the Q-column (and
# B35) literals are written first, in the exact sequence the original
# workbook used, before the remaining numeric/reused-string cells.
# ---------------------------------------------------------------------------

$ws.Range("B35").Value = "FE"
$ws.Range("Q36").Value = "0.7054/0.4058"
$ws.Range("Q42").Value = "0.7055/0.4032"
$ws.Range("Q37").Value = "0.4908/0.4446"
$ws.Range("Q43").Value = "0.4865/0.4259"
$ws.Range("Q24").Value = "0.7333/0.4498"
$ws.Range("Q30").Value = " 0.7316/0.4242"
$ws.Range("Q25").Value = "0.5227/0.4743"
$ws.Range("Q31").Value = "0.5178/0.4511"

# --- Row 24 ---
$ws.Range("B24").Value = "freq"
$ws.Range("D24").Value = 7
$ws.Range("E24").Value = 50
$ws.Range("F24").Value = "z-norm"
$ws.Range("G24").Value = 0.0001
$ws.Range("H24").Value = 200
$ws.Range("I24").Value = "yes"
$ws.Range("J24").Value = 10
$ws.Range("K24").Value = 8
$ws.Range("N24").Value = 662408
$ws.Range("O24").Value = 207014
$ws.Range("P24").Value = 199976
$ws.Range("S24").Value = 0.5135266

# --- Row 25 ---
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = "z-norm"
$ws.Range("G25").Value = 0.00001
$ws.Range("H25").Value = 200
$ws.Range("J25").Value = 10
$ws.Range("K25").Value = 8
$ws.Range("N25").Value = 662408
$ws.Range("O25").Value = 207014
$ws.Range("P25").Value = 199976
$ws.Range("S25").Value = 0.5469856

# --- Row 26 ---
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = "z-norm"
$ws.Range("G26").Value = 0.000001
$ws.Range("H26").Value = 200
$ws.Range("J26").Value = 10
$ws.Range("K26").Value = 8
$ws.Range("N26").Value = 662408
$ws.Range("O26").Value = 207014
$ws.Range("P26").Value = 199976

# --- Row 30 ---
$ws.Range("D30").Value = 8
$ws.Range("E30").Value = 50
$ws.Range("F30").Value = "z-norm"
$ws.Range("G30").Value = 0.0001
$ws.Range("H30").Value = 200
$ws.Range("I30").Value = "yes"
$ws.Range("J30").Value = 10
$ws.Range("K30").Value = 8
$ws.Range("N30").Value = 579594
$ws.Range("O30").Value = 181125
$ws.Range("P30").Value = 174966
$ws.Range("S30").Value = 0.4715259

# --- Row 31 ---
$ws.Range("E31").Value = 50
$ws.Range("F31").Value = "z-norm"
$ws.Range("G31").Value = 0.00001
$ws.Range("H31").Value = 200
$ws.Range("I31").Value = "yes"
$ws.Range("J31").Value = 10
$ws.Range("K31").Value = 8
$ws.Range("N31").Value = 579594
$ws.Range("O31").Value = 181125
$ws.Range("P31").Value = 174966
$ws.Range("S31").Value = 0.5178206

# --- Row 32 ---
$ws.Range("E32").Value = 50
$ws.Range("F32").Value = "z-norm"
$ws.Range("G32").Value = 0.000001
$ws.Range("H32").Value = 200
$ws.Range("I32").Value = "yes"
$ws.Range("J32").Value = 10
$ws.Range("K32").Value = 8
$ws.Range("N32").Value = 579594
$ws.Range("O32").Value = 181125
$ws.Range("P32").Value = 174966

# --- Row 36 ---
$ws.Range("B36").Value = "freq"
$ws.Range("D36").Value = 7
$ws.Range("E36").Value = 50
$ws.Range("F36").Value = "z-norm"
$ws.Range("G36").Value = 0.0001
$ws.Range("H36").Value = 200
$ws.Range("I36").Value = "yes"
$ws.Range("J36").Value = 10
$ws.Range("K36").Value = 8
$ws.Range("N36").Value = 662408
$ws.Range("O36").Value = 207014
$ws.Range("P36").Value = 199976
$ws.Range("S36").Value = 0.44905114

# --- Row 37 ---
$ws.Range("E37").Value = 50
$ws.Range("F37").Value = "z-norm"
$ws.Range("G37").Value = 0.00001
$ws.Range("H37").Value = 200
$ws.Range("I37").Value = "yes"
$ws.Range("J37").Value = 10
$ws.Range("K37").Value = 8
$ws.Range("N37").Value = 662408
$ws.Range("O37").Value = 207014
$ws.Range("P37").Value = 199976
$ws.Range("S37").Value = 0.50065756

# --- Row 38 ---
$ws.Range("E38").Value = 50
$ws.Range("F38").Value = "z-norm"
$ws.Range("G38").Value = 0.000001
$ws.Range("H38").Value = 200
$ws.Range("I38").Value = "yes"
$ws.Range("J38").Value = 10
$ws.Range("K38").Value = 8
$ws.Range("N38").Value = 662408
$ws.Range("O38").Value = 207014
$ws.Range("P38").Value = 199976

# --- Row 42 ---
$ws.Range("D42").Value = 8
$ws.Range("E42").Value = 50
$ws.Range("F42").Value = "z-norm"
$ws.Range("G42").Value = 0.0001
$ws.Range("H42").Value = 200
$ws.Range("I42").Value = "yes"
$ws.Range("J42").Value = 10
$ws.Range("K42").Value = 8
$ws.Range("N42").Value = 579594
$ws.Range("O42").Value = 181125
$ws.Range("P42").Value = 174966
$ws.Range("S42").Value = 0.44833538

# --- Row 43 ---
$ws.Range("E43").Value = 50
$ws.Range("F43").Value = "z-norm"
$ws.Range("G43").Value = 0.00001
$ws.Range("H43").Value = 200
$ws.Range("I43").Value = "yes"
$ws.Range("J43").Value = 10
$ws.Range("K43").Value = 8
$ws.Range("N43").Value = 579594
$ws.Range("O43").Value = 181125
$ws.Range("P43").Value = 174966
$ws.Range("S43").Value = 0.48139343

# --- Row 44 ---
$ws.Range("E44").Value = 50
$ws.Range("F44").Value = "z-norm"
$ws.Range("G44").Value = 0.000001
$ws.Range("H44").Value = 200
$ws.Range("I44").Value = "yes"
$ws.Range("J44").Value = 10
$ws.Range("K44").Value = 8
$ws.Range("N44").Value = 579594
$ws.Range("O44").Value = 181125
$ws.Range("P44").Value = 174966

# ---------------------------------------------------------------------------
# 3) View state: scroll window + active selection
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("K36").Select()
